$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9483",     "Under Review",        ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9484",     "Approved",            ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9485",     "Submitted",           ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9486",     "Disapproved",         ""),
    @("Waiver", "MD", "1915(c)",      "Amendment", "MD-2260.R00.45", "Pending-Approval",    "MD-2260.R00.00"),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9487",     "Pending-Concurrence", ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9488",     "RAI Issued",          ""),
    @("SPA",    "MD", "CHIP SPA",     "",          "MD-25-9489",     "Submitted",           ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9490",     "Submitted",           ""),
    @("Waiver", "MD", "1915(c)",      "Amendment", "MD-2260.R00.46", "Unsubmitted",         "MD-2260.R00.00"),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9491",     "Under Review",        ""),
    @("Waiver", "MD", "1915(b)",      "Initial",   "MD-2279.R00.00", "Terminated",          ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9492",     "Withdrawn",           "")
)

$startRow = 23
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
